# "Added Test Data for UK Market"
# Duplicate the "Netherlands" sheet (placing the copy right after it), rename the
# copy to "UK", and update its test-case id / market name cells.

$wb  = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Netherlands")

# Create a copy of the Netherlands sheet, inserted immediately after it.
# This copy becomes the new active sheet, matching Excel's "Move or Copy -> Create a copy" behavior.
$src.Copy($null, $src)

# The newly created sheet is positioned right after the source sheet.
$uk = $wb.Worksheets.Item($src.Index + 1)
$uk.Name = "UK"

# Update the cells that differ between Netherlands and UK.
# Set B4 before B2 so new shared strings are appended in the same order as the target workbook.
$uk.Range("B4").Value = "NGC-2741/T3392"
$uk.Range("B2").Value = "UK Market"
